$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6.037570333333334, 18.112711, 0.009516155302887736, 0.009516155302887734, 3.0, 1.0, 0.971099, 2.913297, 0.007882574716876797, 0.007882574716876797, 5.863078513129667, 52.767706618167, 0.00007501180519241593, 0.00007501180519241592),
    @(6.037570333333334, 18.112711, 0.009516155302887736, 0.009516155302887734, 3.0, 1.0, 85.826024, 257.478072, 0.6966643430097871, 0.696664343009787, 518.1806563303547, 4663.625906973192, 0.006629566082065386, 0.006629566082065384),
    @(6.037570333333334, 18.112711, 0.009516155302887736, 0.009516155302887734, 3.0, 1.0, 36.24916566666667, 108.747497, 0.294240604502677, 0.294240604502677, 218.8568872371519, 1969.711985134367, 0.002800039288863043, 0.002800039288863042),
    @(6.037570333333334, 18.112711, 0.009516155302887736, 0.009516155302887734, 2.0, 0.6666666666666666, 0.149372, 0.448116, 0.001212477770659141, 0.001212477770659141, 0.9018439558306666, 8.116595602476, 0.00001153812676689148, 0.00001153812676689148),
    @(601.582784, 1804.748352, 0.948188573221463, 0.9481885732214629, 3.0, 1.0, 0.971099, 2.913297, 0.007882574716876797, 0.007882574716876797, 584.196439959616, 5257.767959636543, 0.007474167274106988, 0.007474167274106987),
    @(601.582784, 1804.748352, 0.948188573221463, 0.9481885732214629, 3.0, 1.0, 85.826024, 257.478072, 0.6966643430097871, 0.696664343009787, 51631.45845757081, 464683.1261181373, 0.6605691694127179, 0.6605691694127177),
    @(601.582784, 1804.748352, 0.948188573221463, 0.9481885732214629, 3.0, 1.0, 36.24916566666667, 108.747497, 0.294240604502677, 0.294240604502677, 21806.87399943055, 196261.8659948749, 0.2789955789672141, 0.2789955789672141),
    @(601.582784, 1804.748352, 0.948188573221463, 0.9481885732214629, 2.0, 0.6666666666666666, 0.149372, 0.448116, 0.001212477770659141, 0.001212477770659141, 89.85962361164798, 808.7366125048319, 0.001149657567424031, 0.001149657567424031),
    @(24.52313633333334, 73.56940900000001, 0.03865229901728497, 0.03865229901728497, 3.0, 1.0, 0.971099, 2.913297, 0.007882574716876797, 0.007882574716876797, 23.81439317016367, 214.329538531473, 0.0003046796349828124, 0.0003046796349828123),
    @(24.52313633333334, 73.56940900000001, 0.03865229901728497, 0.03865229901728497, 3.0, 1.0, 85.826024, 257.478072, 0.6966643430097871, 0.696664343009787, 2104.723287499939, 18942.50958749945, 0.02692767850069467, 0.02692767850069467),
    @(24.52313633333334, 73.56940900000001, 0.03865229901728497, 0.03865229901728497, 3.0, 1.0, 36.24916566666667, 108.747497, 0.294240604502677, 0.294240604502677, 888.9432316132528, 8000.489084519275, 0.01137307582826416, 0.01137307582826416),
    @(24.52313633333334, 73.56940900000001, 0.03865229901728497, 0.03865229901728497, 2.0, 0.6666666666666666, 0.149372, 0.448116, 0.001212477770659141, 0.001212477770659141, 3.663069920382667, 32.967629283444, 0.00004686505334332818, 0.00004686505334332817),
    @(2.311301333333333, 6.933904, 0.003642972458364431, 0.003642972458364431, 3.0, 1.0, 0.971099, 2.913297, 0.007882574716876797, 0.007882574716876797, 2.244502413498667, 20.200521721488, 0.00002871600259458198, 0.00002871600259458198),
    @(2.311301333333333, 6.933904, 0.003642972458364431, 0.003642972458364431, 3.0, 1.0, 85.826024, 257.478072, 0.6966643430097871, 0.696664343009787, 198.3698037058987, 1785.328233353088, 0.002537929014309205, 0.002537929014309205),
    @(2.311301333333333, 6.933904, 0.003642972458364431, 0.003642972458364431, 3.0, 1.0, 36.24916566666667, 108.747497, 0.294240604502677, 0.294240604502677, 83.78274493758755, 754.044704438288, 0.001071910418335754, 0.001071910418335753),
    @(2.311301333333333, 6.933904, 0.003642972458364431, 0.003642972458364431, 2.0, 0.6666666666666666, 0.149372, 0.448116, 0.001212477770659141, 0.001212477770659141, 0.3452437027626666, 3.107193324864, 0.000004417023124890354, 0.000004417023124890354)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $colNum = $j + 7
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}